{"js": "// Update Table 2 figures (Std_Error, CI_Lower_95, CI_Upper_95) for the\n// four data rows, per \"updated figure 2 since caption was cut off\".\nconst replacements = [\n  [\"1.44\", \"1.49\"],\n  [\"-4.80\", \"-4.90\"],\n  [\"0.86\", \"0.95\"],\n  [\"1.47\", \"1.40\"],\n  [\"-5.00\", \"-4.85\"],\n  [\"0.78\", \"0.63\"],\n  [\"1.63\", \"1.65\"],\n  [\"-6.11\", \"-6.16\"],\n  [\"0.26\", \"0.32\"],\n  [\"2.29\", \"2.24\"],\n  [\"-5.68\", \"-5.58\"],\n  [\"3.31\", \"3.21\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(`Expected exactly 1 match for \"${oldText}\", found ${results.items.length}`);\n  }\n\n  results.items[0].insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Update Table 2 figures (Std_Error, CI_Lower_95, CI_Upper_95) for the\n# four data rows, per \"updated figure 2 since caption was cut off\".\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n# Row indices are 1-based and include the header row (row 1).\n# Columns: 1=Row, 2=ATT_pp, 3=Std_Error, 4=CI_Lower_95, 5=CI_Upper_95\n$updates = @(\n    @{ Row = 2; Col = 3; Value = \"1.49\" },\n    @{ Row = 2; Col = 4; Value = \"-4.90\" },\n    @{ Row = 2; Col = 5; Value = \"0.95\" },\n    @{ Row = 3; Col = 3; Value = \"1.40\" },\n    @{ Row = 3; Col = 4; Value = \"-4.85\" },\n    @{ Row = 3; Col = 5; Value = \"0.63\" },\n    @{ Row = 4; Col = 3; Value = \"1.65\" },\n    @{ Row = 4; Col = 4; Value = \"-6.16\" },\n    @{ Row = 4; Col = 5; Value = \"0.32\" },\n    @{ Row = 5; Col = 3; Value = \"2.24\" },\n    @{ Row = 5; Col = 4; Value = \"-5.58\" },\n    @{ Row = 5; Col = 5; Value = \"3.21\" }\n)\n\nforeach ($u in $updates) {\n    $cell = $tbl.Cell($u.Row, $u.Col)\n    $cell.Range.Text = $u.Value\n}\n"}
